$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Merge the run-fragmented first body paragraph ("Big data are data
#    sets ...") into a single run. A Find/Replace (even a no-op
#    replace) over the matched range causes same-formatted adjacent
#    runs to be consolidated into one run, matching the target markup.
# ---------------------------------------------------------------------
$old1 = "Big data are data sets which are so large or complex that traditional data processing or learning applications will lead to a poor performance. Some challenges include capturing data, data storage, data analysis, search, sharing, transfer, visualization, querying, updating, etc. Big data causes computational difficulties and intrinsic statistical difficulties due to the data set being of large dimensions. This can cause overfitting, false structures, data isolation, etc. As data grows day by day, exploring different ways of dimension reduction is essential. The aim of my project will be to implement few generic/global fast projections that will reduce the dimensionality of a data set. Thus"
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $old1, 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Restructure the "Extensions" / "Summary of Risks" tail section.
#
#    Before:
#      15 Insufficient knowledge to understand areas in research papers
#      16 Extensions:
#      17 Lack of having a clear idea on the end goal
#      18 (empty paragraph)
#      19 Summary of Risks
#      20 Lack of understanding the mathematics behind the implementation
#      21 Lack of resources to implement projections
#      22 Unforeseen circumstances (... ) [contains the _GoBack bookmark]
#      23 Supervisor becomes unavailable
#
#    After:
#      15 Insufficient knowledge to understand areas in research papers
#         [+ _GoBack bookmark appended at the end]
#      16 Summary of Risks
#      17 Lack of understanding the mathematics behind the implementation
# ---------------------------------------------------------------------

# Locate paragraph 15 ("Insufficient knowledge ...") by its text so this
# does not depend on paragraph numbering staying stable.
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "Insufficient knowledge to understand areas in research papers`r") {
        $targetPara = $d.Paragraphs.Item($i)
        break
    }
}

# Append a throwaway marker character right after the paragraph's text
# (but before its paragraph mark) so we have a safe, non-boundary
# position to anchor the new bookmark on - placing a zero-length
# bookmark exactly at a paragraph-end boundary is ambiguous and lands
# in the following paragraph instead.
$endRange = $targetPara.Range.Duplicate
$endRange.MoveEnd(1, -1)
$endRange.Collapse(0)
$endRange.InsertAfter("@")

$markerRange = $d.Content.Duplicate
$markerRange.Find.Execute("papers@")
$bmAnchor = $d.Range($markerRange.End - 1, $markerRange.End - 1)

# Remove the old "Extensions:" / "Lack of having a clear idea ..." /
# blank paragraph that used to sit between "Insufficient knowledge..."
# and "Summary of Risks" - delete from the bottom up so indices stay
# valid while we work.
$extPara = $null
$blankPara = $null
$lackGoalPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -eq "Extensions:`r") { $extPara = $d.Paragraphs.Item($i) }
    if ($t -eq "Lack of having a clear idea on the end goal`r") { $lackGoalPara = $d.Paragraphs.Item($i) }
}
$blankPara = $lackGoalPara.Next()

$blankPara.Range.Delete()
$lackGoalPara.Range.Delete()
$extPara.Range.Delete()

# Remove "Lack of resources to implement projections", the "Unforeseen
# circumstances (...)" paragraph (which carries the old _GoBack
# bookmark) and "Supervisor becomes unavailable".
$resourcesPara = $null
$unforeseenPara = $null
$supervisorPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -eq "Lack of resources to implement projections`r") { $resourcesPara = $d.Paragraphs.Item($i) }
    if ($t -like "Unforeseen circumstances*") { $unforeseenPara = $d.Paragraphs.Item($i) }
    if ($t -eq "Supervisor becomes unavailable`r") { $supervisorPara = $d.Paragraphs.Item($i) }
}

$supervisorPara.Range.Delete()
$unforeseenPara.Range.Delete()
$resourcesPara.Range.Delete()

# Now add the new _GoBack bookmark right before the temporary "@"
# marker (still safely inside "Insufficient knowledge ..."), then
# remove the marker character.
$d.Bookmarks.Add("_GoBack", $bmAnchor) | Out-Null

$markerRange2 = $d.Content.Duplicate
$markerRange2.Find.Execute("papers@")
$markerChar = $d.Range($markerRange2.End - 1, $markerRange2.End)
$markerChar.Delete()

Write-Output "done"
